# Updated cryptos list on Mon Nov  6 07:06:28 UTC 2023 with GitHub Actions
#
# The "Price" column (D) holds plain text (not numbers) in the source data.
# Some new price strings look numeric (e.g. "246.50", "0.685") and Excel's
# COM layer auto-converts such input to a real number (dropping the text
# formatting / trailing zeros) unless the cell is explicitly formatted as
# Text first. So for any new Price value that would otherwise be parsed as
# a number, we set NumberFormat = "@" immediately before writing it so the
# literal string (including trailing zeros) is preserved. Values that are
# already non-numeric-looking (contain a thousands separator dot, e.g.
# "35.235.67") are written directly since Excel keeps those as text anyway.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $value
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "35.235.67"
$ws.Range("E2").Value = "  -0.84%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.896.41"
$ws.Range("E3").Value = "  -0.32%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.30%  "

# Row 5 - BNB
Set-TextValue "D5" "246.50"
$ws.Range("E5").Value = "  +0.14%  "

# Row 6 - XRP
Set-TextValue "D6" "0.685"
$ws.Range("E6").Value = "  +8.42%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.31%  "

# Row 8 - Solana
Set-TextValue "D8" "40.73"
$ws.Range("E8").Value = "  -3.39%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +2.94%  "

# Row 10 - OKB
Set-TextValue "D10" "52.08"
$ws.Range("E10").Value = "  +7.49%  "

# Row 11 - Dogecoin
Set-TextValue "D11" "0.0720"
$ws.Range("E11").Value = "  +2.35%  "

# Row 12 - TRON
Set-TextValue "D12" "0.0983"
$ws.Range("E12").Value = "  -1.36%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "2.171.07"
$ws.Range("E13").Value = "  -0.38%  "

# Row 14 - Chainlink
Set-TextValue "D14" "12.74"
$ws.Range("E14").Value = "  +2.94%  "

# Row 15 - Polygon
Set-TextValue "D15" "0.708"
$ws.Range("E15").Value = "  +2.64%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "1.899.36"
$ws.Range("E16").Value = "  -0.43%  "

# Row 17 - Polkadot
Set-TextValue "D17" "4.80"
$ws.Range("E17").Value = "  -0.66%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "35.231.59"
$ws.Range("E18").Value = "  -0.81%  "

# Row 19 - Litecoin
Set-TextValue "D19" "71.94"
$ws.Range("E19").Value = "  +0.15%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  +0.82%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "240.50"
$ws.Range("E21").Value = "  -1.18%  "

# Row 22 - Avalanche
Set-TextValue "D22" "12.76"
$ws.Range("E22").Value = "  +2.78%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  -1.74%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  -0.34%  "

# Row 25 - Toncoin
Set-TextValue "D25" "2.32"
$ws.Range("E25").Value = "  +1.49%  "

# Row 26 - PancakeSwap
Set-TextValue "D26" "2.35"
$ws.Range("E26").Value = "  +3.28%  "

# Row 27 - Monero
Set-TextValue "D27" "167.60"
$ws.Range("E27").Value = "  -2.44%  "

# Row 28 - Cosmos
$ws.Range("E28").Value = "  +0.68%  "

# Row 29 - EthereumClassic
Set-TextValue "D29" "18.29"
$ws.Range("E29").Value = "  +1.90%  "

# Row 30 - Stellar
$ws.Range("E30").Value = "  +3.81%  "

# Row 31 - EURNeutrino
$ws.Range("E31").Value = "  +20.02%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +1.61%  "

# Row 33 - Hedera
Set-TextValue "D33" "0.0566"
$ws.Range("E33").Value = "  +0.16%  "

# Row 34 - WEMIXToken
Set-TextValue "D34" "1.89"
$ws.Range("E34").Value = "  +7.94%  "

# Row 35 - BinanceUSD
$ws.Range("E35").Value = "  -0.30%  "

# Row 36 - InternetComputer(DFINITY)
$ws.Range("E36").Value = "  -0.45%  "

# Row 37 - ImmutableX
Set-TextValue "D37" "0.913"
$ws.Range("E37").Value = "  -6.00%  "

# Row 38 - TrustWalletToken
Set-TextValue "D38" "1.50"
$ws.Range("E38").Value = "  +14.92%  "

# Row 39 - LidoDAOToken
$ws.Range("E39").Value = "  -0.15%  "

# Row 40 - InjectiveProtocol
Set-TextValue "D40" "16.43"
$ws.Range("E40").Value = "  +6.05%  "

# Row 41 - ARBITRUM
$ws.Range("E41").Value = "  -1.18%  "

# Row 42 & 43 - VeChain and Kaspa swap places (Kaspa moves to rank 41 / row 42,
# VeChain moves to rank 42 / row 43), with updated values
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D42" "0.0646"
$ws.Range("E42").Value = "  +9.34%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D43" "0.0207"
$ws.Range("E43").Value = "  +0.90%  "

# Row 44 - Aave
Set-TextValue "D44" "90.72"
$ws.Range("E44").Value = "  -0.27%  "

# Row 45 - Maker
$ws.Range("D45").Value = "1.346.46"
$ws.Range("E45").Value = "  -0.28%  "

# Row 46 - RenderToken
Set-TextValue "D46" "2.43"
$ws.Range("E46").Value = "  +4.05%  "

# Row 47 - HuobiToken
$ws.Range("E47").Value = "  +0.00%  "

# Row 48 - MXToken
Set-TextValue "D48" "2.79"
$ws.Range("E48").Value = "  +1.15%  "

# Row 49 - MultiversX
Set-TextValue "D49" "45.80"
$ws.Range("E49").Value = "  -10.28%  "

# Row 50 - Gas
Set-TextValue "D50" "12.08"
$ws.Range("E50").Value = "  -4.35%  "

# Row 51 - FraxShare
$ws.Range("E51").Value = "  -2.70%  "
